# Ajout du fichier de saisie des données de capteurs :
# crée une nouvelle feuille "Àpropos" en fin de classeur contenant
# une date de dépôt (format AAAAMMJJ) en A1.
$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Àpropos"
$newSheet.Range("A1").Value = 20201218
[void]$newSheet.Range("A2").Select()

# Garder la première feuille comme feuille active, comme dans le classeur d'origine.
$wb.Worksheets.Item(1).Activate()
